$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values remain stored as literal text (matching the source
# inlineStr cells) instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "68.186.28"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "3.730.39"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "591.84"
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("D6").Value = "166.10"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").Value = "3.728.57"
$ws.Range("E7").Value = "  -0.64%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -0.63%  "
$ws.Range("D10").Value = "0.160"
$ws.Range("E10").Value = "  -0.15%  "
$ws.Range("E11").Value = "  -0.28%  "
$ws.Range("E12").Value = "  -0.28%  "
$ws.Range("E13").Value = "  -2.03%  "
$ws.Range("D14").Value = "36.10"
$ws.Range("E14").Value = "  +0.30%  "
$ws.Range("D15").Value = "4.359.73"
$ws.Range("E15").Value = "  -0.39%  "
$ws.Range("D16").Value = "3.744.48"
$ws.Range("E16").Value = "  -0.30%  "
$ws.Range("D17").Value = "68.263.20"
$ws.Range("E17").Value = "  +0.55%  "
$ws.Range("D18").Value = "17.88"
$ws.Range("E18").Value = "  -3.20%  "
$ws.Range("E19").Value = "  -0.82%  "
$ws.Range("E20").Value = "  +0.21%  "
$ws.Range("E21").Value = "  +0.76%  "
$ws.Range("D22").Value = "464.91"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").Value = "0.695"
$ws.Range("E23").Value = "  -1.14%  "
$ws.Range("D24").Value = "83.78"
$ws.Range("E24").Value = "  +0.97%  "
$ws.Range("D25").Value = "0.0000146"
$ws.Range("E25").Value = "  +7.59%  "
$ws.Range("E26").Value = "  -0.91%  "
$ws.Range("E27").Value = "  -0.96%  "
$ws.Range("D28").Value = "10.02"
$ws.Range("E28").Value = "  -1.42%  "
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("D30").Value = "3.887.06"
$ws.Range("E30").Value = "  -0.19%  "
$ws.Range("E31").Value = "  -3.85%  "
$ws.Range("E32").Value = "  -1.48%  "
$ws.Range("D33").Value = "29.68"
$ws.Range("E33").Value = "  -0.64%  "
$ws.Range("E34").Value = "  -2.09%  "
$ws.Range("D35").Value = "9.15"
$ws.Range("E35").Value = "  +1.07%  "
$ws.Range("D36").Value = "1.00"
$ws.Range("D37").Value = "3.686.97"
$ws.Range("E37").Value = "  -0.28%  "
$ws.Range("E38").Value = "  -1.84%  "
$ws.Range("D39").Value = "3.36"
$ws.Range("E39").Value = "  -2.13%  "
$ws.Range("D40").Value = "0.137"
$ws.Range("E40").Value = "  +0.30%  "
$ws.Range("D41").Value = "0.994"
$ws.Range("E41").Value = "  -0.46%  "
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("D45").Value = "43.89"
$ws.Range("E45").Value = "  +14.61%  "
$ws.Range("D46").Value = "0.299"
$ws.Range("E46").Value = "  -2.33%  "
$ws.Range("D47").Value = "46.55"
$ws.Range("E47").Value = "  +2.88%  "
$ws.Range("D48").Value = "1.90"
$ws.Range("E48").Value = "  -0.50%  "
$ws.Range("D49").Value = "8.44"
$ws.Range("E49").Value = "  -1.23%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "144.25"
$ws.Range("E50").Value = "  -0.39%  "
$ws.Range("B51").Value = "Bittensor"
$ws.Range("C51").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D51").Value = "387.63"
$ws.Range("E51").Value = "  -2.08%  "
